$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark sitting at the end of the
#    paragraph that ends "...post 2010 slowdown. "
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Expand the single sentence about the Bayes Factor model choice
#    into the fuller, multi-sentence explanation, and re-insert the
#    "_GoBack" bookmark at its new position (right after "...each of
#    these slowdown ", before "models, and how the inclusion...").
# ---------------------------------------------------------------------
$oldRun = "The most likely of this family of alternative models can then be identified, and with each new annual lifetable for the UK and constituent nations the preferred model and changing strength of evidence in support of this model can be updated. A graphical illustration showing the relative likelihood of each of these models, and how the inclusion of each new observation changes the likelihood surface, is shown in Figure R5.1A in the "

$partA = "The most likely of this family of alternative models can then be identified, and with each new annual lifetable for the UK and constituent nations the preferred model and changing strength of evidence in support of this model can be updated."
$partB = " This strength of evidence is expressed as a Bayes Factor, which shows the ratio of the likelihood of two models. In the results presented, a Bayes Factor above 1 indicates more support for a model positing a slowdown from pre-2010 trends, and a ratio below 1 indicates more support for ‘no slowdown’ than ‘slowdown’. "
$partC = "A graphical illustration showing the relative likelihood of each of these "
$partD = "slowdown "
$partE = "models, and how the inclusion of each new observation changes the likelihood surface, is shown in Figure R5.1A in the "

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Execute($oldRun, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$runStart = $findRange.Start
$textBeforeBookmark = $partA + $partB + $partC + $partD
$textAfterBookmark = $partE

$wholeRange = $d.Range($runStart, $findRange.End)
$wholeRange.Text = $textBeforeBookmark + $textAfterBookmark

$bmPos = $runStart + $textBeforeBookmark.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 3) Remove the two empty paragraphs that sit between the paragraph
#    ending "...technical description of the approach. " and the
#    paragraph beginning "For all populations...".
# ---------------------------------------------------------------------
$emptyPara = $d.Paragraphs(4)
$emptyPara.Range.Delete()
$emptyPara2 = $d.Paragraphs(4)
$emptyPara2.Range.Delete()

# ---------------------------------------------------------------------
# 4) Merge the "For all populations..." paragraph's separate runs
#    (split apart by proofErr grammar-check markers around the word
#    "representing") back into a single plain run.
# ---------------------------------------------------------------------
$mergedText = "For all populations except males in Northern Ireland, the addition of the 2018 single year life expectancy data led to sizeable increases in the empirical support for the belief that there has been a slowdown in life expectancy after 2010; this is seen by noting how much higher the bold line, which incorporates the 2018 data, is than the fainter lines representing cumulative data based on shorter series of observations. For most of these populations, the peak of the bold line is to the left of peaks based on earlier series, meaning not only did the 2018 observations increase the strength of evidence supporting belief in a slowdown in life expectancy improvements, but also suggested more severe magnitudes of slowdown than the series excluding this most recent observation had indicated. For the UK as a whole, the addition of the life expectancy data for 2018 suggested an overall slowdown of around 60% was most likely, compared with a most likely magnitude of slowdown of around 50% based on data up to 2017. For each of these populations, what does the Bayes Factor maximise at?"

$popPara = $d.Paragraphs(4)
$popRange = $popPara.Range
$popContent = $d.Range($popRange.Start, $popRange.End - 1)

# The merged text is character-for-character identical to the
# paragraph's existing visible text (only the run/proofErr splitting
# changes), so a plain re-assignment would be treated as a no-op and
# leave the old run/proofErr boundaries in place. Force a genuine
# text diff by appending a sentinel character, then trimming it off,
# so the paragraph is rebuilt as one clean run.
$popContent.Text = $mergedText + "#"
$popPara2 = $d.Paragraphs(4)
$popRange2 = $popPara2.Range
$sentinel = $d.Range($popRange2.End - 2, $popRange2.End - 1)
$sentinel.Text = ""
